# Weekly update for "Macroferia Regional de Talca - Espárragos":
# a new week's record is inserted at row 102 (pushing the existing
# rows 102-114 down to 103-115), and the new row is populated with
# that week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 102; existing rows 102-114 shift to 103-115.
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with this week's data.
$ws.Cells.Item(102, 1).Value = 5
$ws.Cells.Item(102, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(102, 3).Value = "Maule"
$ws.Cells.Item(102, 4).Value = 45211
$ws.Cells.Item(102, 5).Value = 7
$ws.Cells.Item(102, 6).Value = 300000000
$ws.Cells.Item(102, 7).Value = "Espárragos"
$ws.Cells.Item(102, 8).Value = "Sin especificar"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 4000
$ws.Cells.Item(102, 11).Value = 1100
$ws.Cells.Item(102, 12).Value = 1100
$ws.Cells.Item(102, 13).Value = 1100
$ws.Cells.Item(102, 14).Value = "`$/kilo"
$ws.Cells.Item(102, 15).Value = "Provincia de Linares"
$ws.Cells.Item(102, 16).Value = 1100
$ws.Cells.Item(102, 17).Value = 1
$ws.Cells.Item(102, 18).Value = "Hortaliza"
